$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New rows to append to the time-tracking log.
$rows = @(
    @{ Row = 5; Date = 42721; Dauer = "4h"; Was = "Einarbeiten in Asp.net core Techniken" },
    @{ Row = 6; Date = 42722; Dauer = "3h"; Was = "Aufsetzen verschiedener Webprojekte, Austesten der MVC Struktur" },
    @{ Row = 7; Date = 42723; Dauer = "6h"; Was = "Erste View Seiten erstellt" },
    @{ Row = 8; Date = 42724; Dauer = "7h"; Was = "Erste vollständige Verbindung von View -> Controller -> Datenbank hergestellt" }
)

foreach ($r in $rows) {
    $rowIndex = $r.Row

    # Copy the existing date cell's style (A3 uses the date number format)
    # to the new date cell, then overwrite its value.
    $ws.Cells.Item(3, 1).Copy($ws.Cells.Item($rowIndex, 1))
    $ws.Cells.Item($rowIndex, 1).Value = $r.Date

    $ws.Cells.Item($rowIndex, 2).Value = $r.Dauer
    $ws.Cells.Item($rowIndex, 3).Value = $r.Was
}

# Update the selection to reflect the next free row, like in the target file.
$null = $ws.Range("A9").Select()
